# Weekly data refresh: a new daily record is inserted at row 10, pushing the
# existing historical rows (10-116) down by one row (11-117). The workbook's
# dimension grows from A1:R116 to A1:R117.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at position 10; existing rows 10:116 shift to 11:117.
$ws.Rows("10:10").Insert()

# Populate the new row 10 with the latest record. Columns that stay constant
# for this market/product/variety/classification (A, B, C, E, F, G, H, I, N,
# O, Q, R) are copied from the row immediately below (the row that used to
# be row 10), while the date and the volume/price figures (D, J, K, L, M, P)
# carry the new observation's values.
$ws.Range("A10").Value = $ws.Range("A11").Value2
$ws.Range("B10").Value = $ws.Range("B11").Value2
$ws.Range("C10").Value = $ws.Range("C11").Value2
$ws.Range("D10").Value = 45111
$ws.Range("E10").Value = $ws.Range("E11").Value2
$ws.Range("F10").Value = $ws.Range("F11").Value2
$ws.Range("G10").Value = $ws.Range("G11").Value2
$ws.Range("H10").Value = $ws.Range("H11").Value2
$ws.Range("I10").Value = $ws.Range("I11").Value2
$ws.Range("J10").Value = 100
$ws.Range("K10").Value = 19000
$ws.Range("L10").Value = 21000
$ws.Range("M10").Value = 20000
$ws.Range("N10").Value = $ws.Range("N11").Value2
$ws.Range("O10").Value = $ws.Range("O11").Value2
$ws.Range("P10").Value = 1111
$ws.Range("Q10").Value = $ws.Range("Q11").Value2
$ws.Range("R10").Value = $ws.Range("R11").Value2
